# Projeto Clima Tempo - append new sensor readings to the "temperatura" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of temperature / humidity history collected by the (now class-based) logger
$ws.Range("A2").Value = "2024-11-21 21:44:29"
$ws.Range("B2").Value = 20

$ws.Range("A3").Value = "2024-11-21 21:49:33"
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = "Alerta Amarelo, Chuvas Intensas"

$ws.Range("A4").Value = "2024-11-21 21:50:15"
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = "Alerta Amarelo, Chuvas Intensas"

$ws.Range("A5").Value = "2024-11-21 21:51:07"
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = "Alerta Amarelo, Chuvas Intensas"

# Column C now holds the longer "UMIDADE DO AR" alert text, so widen it to fit
$ws.Columns.Item(3).ColumnWidth = 29.45

# Restore the cursor/selection to B12, as left by the author after editing
$ws.Range("B12").Select() | Out-Null
